# Update profit/price figures in the Bahamut Profits sheets (scheduled data refresh)
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 76
$ws.Range("H76").Value = 44003196
$ws.Range("I76").Value = 45836456
$ws.Range("J76").Value = 5001
$ws.Range("K76").Value = 45836456
$ws.Range("L76").Value = 5001
$ws.Range("M76").Value = -5631
# Row 79
$ws.Range("H79").Value = 44003196
$ws.Range("I79").Value = 45836456
$ws.Range("J79").Value = 5001
$ws.Range("K79").Value = 45836456
$ws.Range("L79").Value = 5001
$ws.Range("M79").Value = -7185
# Row 86
$ws.Range("H86").Value = 66670290
$ws.Range("I86").Value = 6001.5
$ws.Range("K86").Value = 6001.5
$ws.Range("M86").Value = -4878.5
# Row 89
$ws.Range("H89").Value = 66670290
$ws.Range("I89").Value = 6001.5
$ws.Range("K89").Value = 30007.5
$ws.Range("M89").Value = -24391.5
# Row 103
$ws.Range("H103").Value = 1127.6666
$ws.Range("I103").Value = 1368.8
$ws.Range("K103").Value = 4106.4
$ws.Range("M103").Value = -3520.4
# Row 129
$ws.Range("H129").Value = 1950411
$ws.Range("I129").Value = 487.33334
$ws.Range("J129").Value = 2850375.8
$ws.Range("K129").Value = 1462.00002
$ws.Range("L129").Value = 8551127.399999999
$ws.Range("M129").Value = 3537.99998
$ws.Range("N129").Value = -8561127.399999999
# Row 135
$ws.Range("H135").Value = 1571.2174
$ws.Range("I135").Value = 862
$ws.Range("J135").Value = 9018
$ws.Range("K135").Value = 7758
$ws.Range("L135").Value = 81162
$ws.Range("M135").Value = -5223
$ws.Range("N135").Value = -86232
# Row 137
$ws.Range("H137").Value = 1091.4615
$ws.Range("I137").Value = 887.5484
$ws.Range("K137").Value = 2662.6452
$ws.Range("M137").Value = -112.6451999999999
# Row 141
$ws.Range("H141").Value = 1968.9117
$ws.Range("I141").Value = 663.37933
$ws.Range("J141").Value = 9541
$ws.Range("K141").Value = 1990.13799
$ws.Range("L141").Value = 28623
$ws.Range("M141").Value = 3189.86201
$ws.Range("N141").Value = -38983

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 3891262.8
$ws.Range("I32").Value = 4235511.5
$ws.Range("J32").Value = 6171.4287
$ws.Range("K32").Value = 4235511.5
$ws.Range("L32").Value = 6171.4287
$ws.Range("M32").Value = -4235224.5
$ws.Range("N32").Value = -6745.4287
# Row 74
$ws.Range("H74").Value = 857.8
$ws.Range("I74").Value = 864.2222
$ws.Range("K74").Value = 864.2222
$ws.Range("M74").Value = 9.777799999999957
# Row 77
$ws.Range("H77").Value = 857.8
$ws.Range("I77").Value = 864.2222
$ws.Range("K77").Value = 4321.111
$ws.Range("M77").Value = 46.88900000000012
# Row 122
$ws.Range("H122").Value = 945.1111
$ws.Range("I122").Value = 934.6667
$ws.Range("J122").Value = 966
$ws.Range("K122").Value = 2804.0001
$ws.Range("L122").Value = 2898
$ws.Range("M122").Value = -354.0001000000002
$ws.Range("N122").Value = -7798
# Row 132
$ws.Range("H132").Value = 1754.7222
$ws.Range("I132").Value = 1006.7692
$ws.Range("J132").Value = 3699.4
$ws.Range("K132").Value = 3020.3076
$ws.Range("L132").Value = 11098.2
$ws.Range("M132").Value = -490.3076000000001
$ws.Range("N132").Value = -16158.2

$ws = $wb.Worksheets.Item("CRP")
# Row 16
$ws.Range("H16").Value = 856.2727
$ws.Range("I16").Value = 703.8
$ws.Range("J16").Value = 983.3333
$ws.Range("K16").Value = 703.8
$ws.Range("L16").Value = 983.3333
$ws.Range("M16").Value = -416.8
$ws.Range("N16").Value = -1557.3333
# Row 31
$ws.Range("H31").Value = 2580.0476
$ws.Range("I31").Value = 2795.1614
$ws.Range("J31").Value = 1973.8182
$ws.Range("K31").Value = 2795.1614
$ws.Range("L31").Value = 1973.8182
$ws.Range("M31").Value = -2500.1614
$ws.Range("N31").Value = -2563.8182
# Row 34
$ws.Range("H34").Value = 2580.0476
$ws.Range("I34").Value = 2795.1614
$ws.Range("J34").Value = 1973.8182
$ws.Range("K34").Value = 2795.1614
$ws.Range("L34").Value = 1973.8182
$ws.Range("M34").Value = -2593.1614
$ws.Range("N34").Value = -2377.8182
# Row 97
$ws.Range("H97").Value = 0
$ws.Range("J97").Value = 0
$ws.Range("L97").Value = 0
$ws.Range("N97").ClearContents()
# Row 113
$ws.Range("H113").Value = 856.2727
$ws.Range("I113").Value = 703.8
$ws.Range("J113").Value = 983.3333
$ws.Range("K113").Value = 703.8
$ws.Range("L113").Value = 983.3333
$ws.Range("M113").Value = 1466.2
$ws.Range("N113").Value = -5323.3333
# Row 132
$ws.Range("H132").Value = 3626
$ws.Range("I132").Value = 1940.5
$ws.Range("J132").Value = 4749.6665
$ws.Range("K132").Value = 5821.5
$ws.Range("L132").Value = 14248.9995
$ws.Range("M132").Value = -3291.5
$ws.Range("N132").Value = -19308.9995
# Row 134
$ws.Range("H134").Value = 6327.04
$ws.Range("I134").Value = 4972.6875
$ws.Range("J134").Value = 8734.777
$ws.Range("K134").Value = 14918.0625
$ws.Range("L134").Value = 26204.331
$ws.Range("M134").Value = -12383.0625
$ws.Range("N134").Value = -31274.331

$ws = $wb.Worksheets.Item("CUL")
# Row 107
$ws.Range("H107").Value = 371153
$ws.Range("I107").Value = 891.3333
$ws.Range("J107").Value = 864835.25
$ws.Range("K107").Value = 2673.9999
$ws.Range("L107").Value = 2594505.75
$ws.Range("M107").Value = -753.9998999999998
$ws.Range("N107").Value = -2598345.75

$ws = $wb.Worksheets.Item("GSM")
# Row 58
$ws.Range("H58").Value = 8970.5
$ws.Range("I58").Value = 8941
$ws.Range("K58").Value = 8941
$ws.Range("M58").Value = -8664
# Row 126
$ws.Range("H126").Value = 3811.25
$ws.Range("I126").Value = 5104
$ws.Range("J126").Value = 3035.6
$ws.Range("K126").Value = 15312
$ws.Range("L126").Value = 9106.799999999999
$ws.Range("M126").Value = -12842
$ws.Range("N126").Value = -14046.8
# Row 132
$ws.Range("H132").Value = 6399.7144
$ws.Range("I132").Value = 7700
$ws.Range("J132").Value = 4666
$ws.Range("K132").Value = 23100
$ws.Range("L132").Value = 13998
$ws.Range("M132").Value = -20570
$ws.Range("N132").Value = -19058

$ws = $wb.Worksheets.Item("LTW")
# Row 82
$ws.Range("H82").Value = 3205.4119
$ws.Range("I82").Value = 3074.9
$ws.Range("J82").Value = 3391.8572
$ws.Range("K82").Value = 3074.9
$ws.Range("L82").Value = 3391.8572
$ws.Range("M82").Value = -2713.9
$ws.Range("N82").Value = -4113.8572
# Row 85
$ws.Range("H85").Value = 3205.4119
$ws.Range("I85").Value = 3074.9
$ws.Range("J85").Value = 3391.8572
$ws.Range("K85").Value = 3074.9
$ws.Range("L85").Value = 3391.8572
$ws.Range("M85").Value = -1826.9
$ws.Range("N85").Value = -5887.8572
# Row 100
$ws.Range("H100").Value = 8548778
$ws.Range("I100").Value = 11112791
$ws.Range("K100").Value = 11112791
$ws.Range("M100").Value = -11112250
# Row 136
$ws.Range("H136").Value = 2710.5151
$ws.Range("I136").Value = 1756.6818
$ws.Range("J136").Value = 4618.1816
$ws.Range("K136").Value = 5270.0454
$ws.Range("L136").Value = 13854.5448
$ws.Range("M136").Value = -2720.0454
$ws.Range("N136").Value = -18954.5448

$ws = $wb.Worksheets.Item("WVR")
# Row 107
$ws.Range("H107").Value = 707.6
$ws.Range("I107").Value = 740.61536
$ws.Range("K107").Value = 2221.84608
$ws.Range("M107").Value = -301.8460800000003
# Row 132
$ws.Range("H132").Value = 2132.8572
$ws.Range("I132").Value = 1690.8948
$ws.Range("J132").Value = 3065.889
$ws.Range("K132").Value = 5072.6844
$ws.Range("L132").Value = 9197.667000000001
$ws.Range("M132").Value = -2542.6844
$ws.Range("N132").Value = -14257.667

Write-Host "Updated 31 rows across 7 sheets"
